$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (pushing current row 8 and below down by one)
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new key/value pair
$ws.Range("A8").Value = "new"
$ws.Range("B8").Value = "* NEW *"

# Update the view: scroll back to top (removes the stale topLeftCell) and select B8
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B8").Select()
